# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 3-21
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 76
$wsExhibit.Range("F4").Value = 1511
$wsExhibit.Range("F5").Value = 581
$wsExhibit.Range("F6").Value = 1070
$wsExhibit.Range("F7").Value = 11085
$wsExhibit.Range("F13").Value = 12235
$wsExhibit.Range("F14").Value = 12780
$wsExhibit.Range("F21").Value = 46

# Sheet "全部类型" (all types) - same events, rows shifted by +1
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 76
$wsAll.Range("F5").Value = 1511
$wsAll.Range("F6").Value = 581
$wsAll.Range("F7").Value = 1070
$wsAll.Range("F8").Value = 11085
$wsAll.Range("F14").Value = 12235
$wsAll.Range("F15").Value = 12780
$wsAll.Range("F22").Value = 46
